$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to make room for the "Match ID" column.
# (This shifts every existing column A..AC one position to the right, B..AD,
#  and Excel automatically grows the sheet dimension / merged-cell refs.)
[void]$ws.Columns("A:A").Insert()

# Header label for the newly inserted column, in the header-labels row (row 3).
$ws.Range("A3").Value = "Match ID"

# Populate the Match ID value (28) for every data row, including the hidden
# totals row 20.
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 28
}

# The header cell and the per-player Match ID values (rows 3-19) pick up a
# bold, border-less style (a new cellXfs entry distinct from the existing
# bordered header style used by the neighbouring "Player ID" column).
$ws.Range("A3:A19").Font.Bold = $true

# Row 20 is hidden and excel auto-sizes hidden rows whenever their content
# changes; re-running AutoFit clears that stray explicit row height so the
# row definition stays exactly as it was (just "hidden").
[void]$ws.Rows(20).AutoFit()

# Match the author's new selection: the whole new Match ID column of visible
# data rows, anchored at its first cell.
[void]$ws.Range("A3:A19").Select()
